$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1").Formula = "=TODAY()"

# Excel auto-sizes the new column to fit the displayed date ("best fit").
$col = $ws.Range("Q1").EntireColumn
$col.AutoFit()
$col.ColumnWidth = 12 - 5/6

# Scroll the view so the new column is visible, and land the selection on it
# (mirrors what Excel does after you type a value in a cell off to the right).
$excel.ActiveWindow.ScrollColumn = 13
$ws.Range("Q1").Select() | Out-Null
